# Insert a new data row at row 424 (pushes existing rows 424:485 down to 425:486)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(424).Insert()

$ws.Range("A424").Value = 10
$ws.Range("B424").Value = "Vega Modelo de Temuco"
$ws.Range("C424").Value = "La Araucanía"
$ws.Range("D424").Value = 45154
$ws.Range("E424").Value = 9
$ws.Range("F424").Value = "Fruta"
$ws.Range("G424").Value = 100102
$ws.Range("H424").Value = "Cítricos"
$ws.Range("I424").Value = 100102006
$ws.Range("J424").Value = "Pomelo"
$ws.Range("K424").Value = "Start Ruby"
$ws.Range("L424").Value = "Primera"
$ws.Range("M424").Value = 45
$ws.Range("N424").Value = 15000
$ws.Range("O424").Value = 15000
$ws.Range("P424").Value = 15000
$ws.Range("Q424").Value = "$/bandeja 15 kilos granel"
$ws.Range("R424").Value = "Región de O'Higgins"
$ws.Range("S424").Value = 1000
$ws.Range("T424").Value = 15

$ws.Range("D424").NumberFormat = "YYYY-MM-DD HH:MM:SS"
